$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 259
$ws.Range("F4").Value = 274
$ws.Range("F5").Value = 2875
$ws.Range("F8").Value = 2232
$ws.Range("F9").Value = 1387
$ws.Range("F13").Value = 2559
$ws.Range("F14").Value = 86
$ws.Range("F15").Value = 1374
$ws.Range("F16").Value = 4761
$ws.Range("F18").Value = 5242
$ws.Range("F19").Value = 1801
$ws.Range("F20").Value = 2913
$ws.Range("F21").Value = 3319
$ws.Range("F22").Value = 181
$ws.Range("F23").Value = 1589
$ws.Range("F24").Value = 265
$ws.Range("F26").Value = 117
$ws.Range("F27").Value = 307
$ws.Range("F29").Value = 1983
$ws.Range("F30").Value = 122
$ws.Range("F31").Value = 292
$ws.Range("F32").Value = 750
$ws.Range("F33").Value = 160
$ws.Range("F34").Value = 352
$ws.Range("F35").Value = 429

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 101
$ws.Range("F10").Value = 24
$ws.Range("F13").Value = 27

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 101
$ws.Range("F7").Value = 259
$ws.Range("F9").Value = 274
$ws.Range("F10").Value = 2875
$ws.Range("F12").Value = 2232
$ws.Range("F13").Value = 1387
$ws.Range("F19").Value = 24
$ws.Range("F20").Value = 2559
$ws.Range("F21").Value = 1374
$ws.Range("F24").Value = 27
$ws.Range("F25").Value = 4761
$ws.Range("F27").Value = 5242
$ws.Range("F28").Value = 1801
$ws.Range("F29").Value = 2913
$ws.Range("F30").Value = 3319
$ws.Range("F32").Value = 181
$ws.Range("F35").Value = 1589
$ws.Range("F37").Value = 265
$ws.Range("F39").Value = 117
$ws.Range("F40").Value = 307
$ws.Range("F43").Value = 1983
$ws.Range("F44").Value = 122
$ws.Range("F45").Value = 292
$ws.Range("F46").Value = 750
$ws.Range("F47").Value = 160
$ws.Range("F48").Value = 352
$ws.Range("F49").Value = 429
